$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conventionals")
$rng = $ws.Range("H1")
$rng.Value = "strike_price"
$rng.Font.Bold = $true
$rng.Font.Name = "Calibri"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
$rng.Borders.Item(7).LineStyle = 1
$rng.Borders.Item(7).Weight = 2
$rng.Borders.Item(7).ColorIndex = -4105
$rng.Borders.Item(10).LineStyle = 1
$rng.Borders.Item(10).Weight = 2
$rng.Borders.Item(10).ColorIndex = -4105
